# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "last row" (row 12) had the short-date style (YYYY-MM-DD).
# Now that a new row is appended, row 12 reverts to the standard
# date-time style used by every other non-final row, and the new row 13
# becomes the final row carrying the short-date style.
$ws.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data.
$ws.Range("A13").Value = 45962
$ws.Range("A13").NumberFormat = "YYYY-MM-DD"
$ws.Range("B13").Value = 24
$ws.Range("C13").Value = 35
$ws.Range("D13").Value = 27
